# The presentation's two themes (theme1.xml used by the slide master,
# theme2.xml used by the notes master) swap their colour schemes:
#   - theme1.xml ("Integral" colours)      -> becomes the "Office" colours
#   - theme2.xml ("Office Theme" colours)  -> becomes the "Integral" colours
# Font scheme / format scheme are identical between the two themes already,
# so only the 12 theme colour scheme slots need to be rewritten on each
# master's ColorScheme (slot order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). RGB values use the standard Windows RGB() packing
# (R + G*256 + B*65536).

$p = $ppt.ActivePresentation

# ---- Slide master (ppt/theme/theme1.xml) becomes the "Office" palette ----
$masterScheme = $p.SlideMaster.ColorScheme
$masterScheme.Colors(1).RGB  = 0         # dk1      000000
$masterScheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$masterScheme.Colors(3).RGB  = 6968388   # dk2      44546A
$masterScheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$masterScheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$masterScheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$masterScheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$masterScheme.Colors(8).RGB  = 49407     # accent4  FFC000
$masterScheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$masterScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$masterScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$masterScheme.Colors(12).RGB = 7491477   # folHlink 954F72

# ---- Notes master (ppt/theme/theme2.xml) becomes the "Integral" palette ----
$notesScheme = $p.NotesMaster.ColorScheme
$notesScheme.Colors(1).RGB  = 0         # dk1      000000
$notesScheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$notesScheme.Colors(3).RGB  = 5332805   # dk2      455F51
$notesScheme.Colors(4).RGB  = 13754083  # lt2      E3DED1
$notesScheme.Colors(5).RGB  = 3722137   # accent1  99CB38
$notesScheme.Colors(6).RGB  = 3646819   # accent2  63A537
$notesScheme.Colors(7).RGB  = 2412774   # accent3  E6D024
$notesScheme.Colors(8).RGB  = 38860     # accent4  CC9700
$notesScheme.Colors(9).RGB  = 13611854  # accent5  4EB3CF
$notesScheme.Colors(10).RGB = 10915127  # accent6  378DA6
$notesScheme.Colors(11).RGB = 2465643   # hlink    6B9F25
$notesScheme.Colors(12).RGB = 158642    # folHlink B26B02
